# Logboek update: add the "PeopleController" debugging entry for 28/12/2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New logboek row (row 14): task description, date, and hours for both students.
$ws.Range("A14").Value = "Debuggen van de PeopleController (afstellen met nieuwe back-end en db)"
$ws.Range("B14").Value = "28/12/2024"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2

# The date column in this sheet is formatted/right-aligned text (see row 13),
# so copy that cell's format onto the new date cell instead of leaving the
# default numeric style behind.
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# Row grows to fit the wrapped task text, same as the other populated rows.
$ws.Rows.Item(14).RowHeight = 58.2

# Reflect where the user ended up after entering the data.
$ws.Range("E9").Select() | Out-Null
